$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Work from the end of the document backwards so earlier paragraph
# --- indices / offsets stay valid while we edit. ---

# 8th paragraph: used to just hold the "_GoBack" bookmark. That bookmark
# moves into paragraph 5 (see below), so this paragraph becomes plain and
# three more blank paragraphs are added after it.
$p8 = $d.Paragraphs(8)
$p8.Range.InsertXML($pkgOpen + '<w:body><w:p/><w:p/><w:p/><w:p/></w:body>' + $pkgClose)

# 7th paragraph: "Table of contents:" becomes a single tab character.
$p7 = $d.Paragraphs(7)
$p7.Range.InsertXML($pkgOpen + '<w:body><w:p><w:r><w:tab/></w:r></w:p></w:body>' + $pkgClose)

# 5th paragraph (centered title): "User Manual for compiler " becomes
# "User Manual for " + "C" + [bookmark _GoBack] + "ompiler ".
$p5 = $d.Paragraphs(5)
$p5body = '<w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">User Manual for </w:t></w:r>' + `
    '<w:r><w:t>C</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t xml:space="preserve">ompiler </w:t></w:r>' + `
    '</w:p></w:body>'
$p5.Range.InsertXML($pkgOpen + $p5body + $pkgClose)
